$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10, pushing the existing rows 10-21 down to 11-22.
$ws.Rows("10:10").Insert()

# Populate the newly inserted row 10 with the new weekly record.
$ws.Range("A10").Value = 11
$ws.Range("B10").Value = "Vega Monumental Concepción"
$ws.Range("C10").Value = "Bíobío"
$ws.Range("D10").Value = 44589
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 100112022
$ws.Range("G10").Value = "Arveja Verde"
$ws.Range("H10").Value = "Perfection"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 160
$ws.Range("K10").Value = 22000
$ws.Range("L10").Value = 23000
$ws.Range("M10").Value = 22500
$ws.Range("N10").Value = "$/malla 25 kilos"
$ws.Range("O10").Value = "Carahue"
$ws.Range("P10").Value = 900
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
